$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Apply the built-in "Comma" cell style (0 decimals) to G26:H26
$ws.Range("G26:H26").Style = "Comma"
$ws.Range("G26:H26").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

# Update H26's value (consumption moved to stock)
$ws.Range("H26").Value = 7500000

# Update sheet view state
$ws.Application.ActiveWindow.DisplayRightToLeft = $true
$ws.Range("H26").Select()
